# Updated cryptos list on Tue May  7 03:11:46 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (D) and Volume(1h) (E) updates, keyed by row number ---
# Each entry: row, newPrice (or $null to leave unchanged), newVolume (or $null to leave unchanged)
# Price values that Excel would otherwise auto-parse as a genuine number (single
# decimal point, plain digits) are written with a leading apostrophe so they are
# kept as literal text -- exactly like the "Price" column already is for every
# other row (e.g. "63.680.78") instead of being silently turned into a Number.
$updates = @(
    @(2,  "63.680.78", "  -0.65%  "),
    @(3,  "3.085.70",  "  -1.99%  "),
    @(4,  $null,       "  -0.29%  "),
    @(5,  "'590.94",   "  -0.08%  "),
    @(6,  "'156.40",   "  +6.80%  "),
    @(7,  $null,       "  -0.10%  "),
    @(8,  $null,       "  +2.45%  "),
    @(9,  "3.084.38",  "  -1.69%  "),
    @(10, "'0.156",    "  -3.63%  "),
    @(11, $null,       "  -1.29%  "),
    @(12, $null,       "  -0.58%  "),
    @(13, "'37.27",    "  +0.04%  "),
    @(14, $null,       "  -3.10%  "),
    @(15, "3.598.52",  "  -1.97%  "),
    @(16, $null,       "  -1.69%  "),
    @(17, "'7.16",     "  -1.30%  "),
    @(18, "63.657.89", "  -0.38%  "),
    @(19, "3.085.91",  "  -1.71%  "),
    @(20, "'478.73",   "  +2.26%  "),
    @(21, "'14.51",    "  +1.15%  "),
    @(22, $null,       "  -2.64%  "),
    @(23, $null,       "  -0.33%  "),
    @(24, "'2.42",     "  +1.64%  "),
    @(25, "'81.20",    "  +0.27%  "),
    @(26, $null,       "  -2.02%  "),
    @(27, "'10.32",    "  +5.09%  "),
    @(28, "'0.998",    "  -0.15%  "),
    @(29, "'7.47",     "  +1.40%  "),
    @(30, "'2.67",     "  -1.62%  "),
    @(31, $null,       "  -0.17%  "),
    @(32, $null,       "  -2.11%  "),
    @(33, $null,       "  +0.23%  "),
    @(34, "'27.27",    "  -1.43%  "),
    @(35, "0.0₃0851",  "  -1.81%  "),
    @(36, "'3.49",     "  +7.72%  "),
    @(37, $null,       "  -1.18%  "),
    @(39, $null,       "  -2.65%  "),
    @(40, $null,       "  -0.28%  "),
    @(41, $null,       "  -1.09%  "),
    @(42, "'446.34",   "  -3.18%  "),
    @(47, "2.815.21",  "  -2.93%  "),
    @(48, "'131.71",   "  -0.60%  "),
    @(49, "'25.43",    "  +5.95%  "),
    @(50, "'1.00",     "  +0.06%  "),
    @(51, $null,       "  +1.25%  ")
)

foreach ($u in $updates) {
    $row = $u[0]
    $price = $u[1]
    $volume = $u[2]
    if ($null -ne $price) {
        $ws.Cells.Item($row, 4).Value = $price
    }
    if ($null -ne $volume) {
        $ws.Cells.Item($row, 5).Value = $volume
    }
}

# --- Rows 43-46 got reshuffled (coins swapped in pairs, with new price data) ---
$ws.Cells.Item(43, 2).Value = "TheGraph"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Cells.Item(43, 4).Value = "'0.288"
$ws.Cells.Item(43, 5).Value = "  -1.69%  "

$ws.Cells.Item(44, 2).Value = "Arweave"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Cells.Item(44, 4).Value = "'41.28"
$ws.Cells.Item(44, 5).Value = "  +3.11%  "

$ws.Cells.Item(45, 2).Value = "Kaspa"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(45, 4).Value = "'0.112"
$ws.Cells.Item(45, 5).Value = "  +4.24%  "

$ws.Cells.Item(46, 2).Value = "VeChain"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(46, 4).Value = "'0.0362"
$ws.Cells.Item(46, 5).Value = "  -2.80%  "
